$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.455.57"
$ws.Range("E2").Value = "  +1.07%  "

$ws.Range("D3").Value = "2.300.71"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.47"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.44%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.27"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("E7").Value = "  +0.71%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.93"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.75%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0908"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.63%  "

$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.964"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.33"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("D16").Value = "2.650.24"
$ws.Range("E16").Value = "  -0.08%  "

$ws.Range("D17").Value = "2.308.09"
$ws.Range("E17").Value = "  +0.49%  "

$ws.Range("D18").Value = "42.433.70"
$ws.Range("E18").Value = "  +0.75%  "

$ws.Range("E19").Value = "  -2.18%  "

$ws.Range("E20").Value = "  +0.93%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.34"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.50%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "276.60"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.65%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.52"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.62%  "

$ws.Range("E24").Value = "  +19.85%  "

$ws.Range("E25").Value = "  -0.94%  "

$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.82"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.77"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "165.27"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "

$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.61"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("E32").Value = "  -2.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.87"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("E34").Value = "  +3.96%  "

$ws.Range("E35").Value = "  -0.73%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.55"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -12.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0367"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.44%  "

$ws.Range("E38").Value = "  +1.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.73"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.53%  "

$ws.Range("E40").Value = "  -0.63%  "

$ws.Range("E41").Value = "  +2.29%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.61"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.97%  "

$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.53"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.71%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "81.58"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.86%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.05"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.85"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.37%  "

$ws.Range("E49").Value = "  -1.24%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.18"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.55%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "1.590.18"
$ws.Range("E51").Value = "  +2.01%  "
